$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column E into column F for all rows so the new
# column inherits the same header/data styles (s=1 header, s=2 data rows)
$ws.Range("E1:E12").Copy()
$ws.Range("F1:F12").PasteSpecial(-4122)

# Add the new "Ảnh" (Image) column with the habitat image file names
$ws.Range("F1").Value = "Ảnh"
$ws.Range("F2").Value = "BaoHabitat.jpg"
$ws.Range("F3").Value = "ChimCongHabitat.jpg"
$ws.Range("F4").Value = "VoiHabitat.jpg"
$ws.Range("F5").Value = "VoiHabitat.jpg"
$ws.Range("F6").Value = "CaSauHabitat.jpg"
$ws.Range("F7").Value = "SuTuHabitat.jpg"
$ws.Range("F8").Value = "KhiHabitat.jpg"
$ws.Range("F9").Value = "SuTuHabitat.jpg"
$ws.Range("F10").Value = "TeGiacHabitat.jpg"
$ws.Range("F11").Value = "HongHacHabiat.jpg"
$ws.Range("F12").Value = "SoiHabitat.jpeg"

# Size column F similarly to the other bestFit columns
$colF = $ws.Range("F1").EntireColumn
$colF.ColumnWidth = 22.8333333
